# Add a new trailing worksheet "ODI Batting Extra" (4th sheet) with a
# header row + one data row, matching the existing header style used on
# the other sheets ("ODI Bowling" header formatting == style index 1:
# bold font, thin border, centered/top aligned).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Reuse the header formatting already present on another sheet instead of
# fabricating a brand-new style.
$srcWs = $wb.Worksheets.Item("ODI Bowling")
$srcWs.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Data row - force text storage (leading apostrophe) so "3861" stays a
# string rather than being coerced to a number, then drop back to the
# default "Normal" style so no extra formatting is left behind.
$ws.Range("A2").Value = "'3861"
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NO"
$ws.Range("A2:E2").Style = "Normal"
